$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price/Volume columns so numeric-looking
# strings (e.g. "597.09") are not auto-converted to numbers by Excel.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.301.66'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.36%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.517.38'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.95%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.09'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.65%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.30'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.45%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.595'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.41%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.134'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +7.31%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.29'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.27%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.437'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.43%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.126.17'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.97%  '

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.28%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.15'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.88%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '67.225.39'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.26%  '

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.26%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.523.36'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.73%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.36'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.54%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.26'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.64%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '395.66'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.39%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.01'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.33%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.22'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.37%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.541'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.16%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.07%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000123'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.64%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.26'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.98%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.183'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.85%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.997'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.27%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.30'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.97%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.47'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.83%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.07'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.56%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '23.95'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.16%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.43'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.57%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.68'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +5.01%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '163.34'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.43%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.896'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.25%  '

# Row 37
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.19'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +8.88%  '

# Row 38
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.92'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.47%  '

# Row 39
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0754'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.16%  '

# Row 40
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.70'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.48%  '

# Row 41
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.67'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.96%  '

# Row 42
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.41'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.76%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.65'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +5.77%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.840.20'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.36%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.94'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.16%  '

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.07%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '340.65'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.24%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.09'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.60%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '34.33'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.86%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.52'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.09%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.854'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.22%  '
